$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 937.125
$ws.Range("I19").Value = 933.2
$ws.Range("J19").Value = 943.6667
$ws.Range("K19").Value = 933.2
$ws.Range("L19").Value = 943.6667
$ws.Range("M19").Value = -758.2
$ws.Range("N19").Value = -1293.6667

$ws.Range("H40").Value = 4999.5
$ws.Range("I40").Value = 4499.5
$ws.Range("K40").Value = 4499.5
$ws.Range("M40").Value = -4324.5

$ws.Range("H92").Value = 4599.5557
$ws.Range("I92").Value = 4542.2856
$ws.Range("J92").Value = 4800
$ws.Range("K92").Value = 4542.2856
$ws.Range("L92").Value = 4800
$ws.Range("M92").Value = -3294.2856
$ws.Range("N92").Value = -7296

$ws.Range("H116").Value = 8896.714
$ws.Range("I116").Value = 12692
$ws.Range("J116").Value = 3836.3333
$ws.Range("K116").Value = 12692
$ws.Range("L116").Value = 3836.3333
$ws.Range("M116").Value = -9250
$ws.Range("N116").Value = -10720.3333

$ws.Range("H131").Value = 3835.875
$ws.Range("I131").Value = 2255.2856
$ws.Range("K131").Value = 6765.8568
$ws.Range("M131").Value = -1725.8568

$ws.Range("H132").Value = 3705.1667
$ws.Range("I132").Value = 3546.025
$ws.Range("K132").Value = 10638.075
$ws.Range("M132").Value = -8108.075000000001

$ws.Range("H137").Value = 1645.8529
$ws.Range("I137").Value = 1475.5
$ws.Range("J137").Value = 2440.8333
$ws.Range("K137").Value = 4426.5
$ws.Range("L137").Value = 7322.499899999999
$ws.Range("M137").Value = -1876.5
$ws.Range("N137").Value = -12422.4999

$ws.Range("H141").Value = 9420.200000000001
$ws.Range("I141").Value = 7946.4614
$ws.Range("K141").Value = 23839.3842
$ws.Range("M141").Value = -18659.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2771.7273
$ws.Range("I32").Value = 1558.4902
$ws.Range("J32").Value = 18240.5
$ws.Range("K32").Value = 1558.4902
$ws.Range("L32").Value = 18240.5
$ws.Range("M32").Value = -1271.4902
$ws.Range("N32").Value = -18814.5

$ws.Range("H45").Value = 3677.3333
$ws.Range("J45").Value = 3128.1428
$ws.Range("L45").Value = 3128.1428
$ws.Range("N45").Value = -3882.1428

$ws.Range("H74").Value = 1433.8302
$ws.Range("I74").Value = 631.4516
$ws.Range("K74").Value = 631.4516
$ws.Range("M74").Value = 242.5484

$ws.Range("H77").Value = 1433.8302
$ws.Range("I77").Value = 631.4516
$ws.Range("K77").Value = 3157.258
$ws.Range("M77").Value = 1210.742

$ws.Range("H103").Value = 59444
$ws.Range("J103").Value = 59444
$ws.Range("L103").Value = 59444
$ws.Range("N103").Value = -61788

$ws.Range("H139").Value = 55344.965
$ws.Range("J139").Value = 55344.965
$ws.Range("L139").Value = 55344.965
$ws.Range("N139").Value = -65624.965

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 5902
$ws.Range("J28").Value = 5902
$ws.Range("L28").Value = 5902
$ws.Range("N28").Value = -6490

$ws.Range("H103").Value = 29901.75
$ws.Range("J103").Value = 29901.75
$ws.Range("L103").Value = 29901.75
$ws.Range("N103").Value = -32245.75

$ws.Range("H106").Value = 17859.092
$ws.Range("J106").Value = 17859.092
$ws.Range("L106").Value = 17859.092
$ws.Range("N106").Value = -20383.092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3445
$ws.Range("I31").Value = 3900
$ws.Range("K31").Value = 3900
$ws.Range("M31").Value = -3605

$ws.Range("H34").Value = 3445
$ws.Range("I34").Value = 3900
$ws.Range("K34").Value = 3900
$ws.Range("M34").Value = -3698

$ws.Range("H105").Value = 1262.7037
$ws.Range("I105").Value = 1243.2609
$ws.Range("J105").Value = 1374.5
$ws.Range("K105").Value = 1243.2609
$ws.Range("L105").Value = 1374.5
$ws.Range("M105").Value = 503.7391
$ws.Range("N105").Value = -4868.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 1639.4
$ws.Range("I15").Value = 56.285713
$ws.Range("J15").Value = 5333.3335
$ws.Range("K15").Value = 168.857139
$ws.Range("L15").Value = 16000.0005
$ws.Range("M15").Value = -28.85713900000002
$ws.Range("N15").Value = -16280.0005

$ws.Range("H34").Value = 3642.8333
$ws.Range("I34").Value = 102.57143
$ws.Range("J34").Value = 8599.200000000001
$ws.Range("K34").Value = 307.71429
$ws.Range("L34").Value = 25797.6
$ws.Range("M34").Value = -223.71429
$ws.Range("N34").Value = -25965.6

$ws.Range("H103").Value = 1313
$ws.Range("I103").Value = 1313
$ws.Range("K103").Value = 3939
$ws.Range("M103").Value = -3060

$ws.Range("H131").Value = 2129
$ws.Range("J131").Value = 2501.3333
$ws.Range("L131").Value = 7503.999899999999
$ws.Range("N131").Value = -17583.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

$ws.Range("H43").Value = 2617
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null

$ws.Range("H70").Value = 9866.177
$ws.Range("I70").Value = 11429.728
$ws.Range("J70").Value = 6999.6665
$ws.Range("K70").Value = 11429.728
$ws.Range("L70").Value = 6999.6665
$ws.Range("M70").Value = -11159.728
$ws.Range("N70").Value = -7539.6665

$ws.Range("H73").Value = 9866.177
$ws.Range("I73").Value = 11429.728
$ws.Range("J73").Value = 6999.6665
$ws.Range("K73").Value = 11429.728
$ws.Range("L73").Value = 6999.6665
$ws.Range("M73").Value = -10493.728
$ws.Range("N73").Value = -8871.666499999999

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8971.714
$ws.Range("J7").Value = 8979.6
$ws.Range("L7").Value = 8979.6
$ws.Range("N7").Value = -9203.6

$ws.Range("H61").Value = 1844.1177
$ws.Range("I61").Value = 1679.6666
$ws.Range("J61").Value = 2238.8
$ws.Range("K61").Value = 1679.6666
$ws.Range("L61").Value = 2238.8
$ws.Range("M61").Value = -1477.6666
$ws.Range("N61").Value = -2642.8

$ws.Range("H69").Value = 220000
$ws.Range("I69").Value = 220000
$ws.Range("K69").Value = 220000
$ws.Range("M69").Value = -219189

$ws.Range("H72").Value = 220000
$ws.Range("I72").Value = 220000
$ws.Range("K72").Value = 660000
$ws.Range("M72").Value = -655944

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null

$ws.Range("H82").Value = 1197.0834
$ws.Range("I82").Value = 499.25
$ws.Range("K82").Value = 499.25
$ws.Range("M82").Value = -138.25

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null

$ws.Range("H85").Value = 1197.0834
$ws.Range("I85").Value = 499.25
$ws.Range("K85").Value = 499.25
$ws.Range("M85").Value = 748.75

$ws.Range("H97").Value = 75000
$ws.Range("J97").Value = 75000
$ws.Range("L97").Value = 75000
$ws.Range("N97").Value = -76982

$ws.Range("H106").Value = 19833
$ws.Range("J106").Value = 19833
$ws.Range("L106").Value = 19833
$ws.Range("N106").Value = -22357

$ws.Range("H113").Value = 1844.1177
$ws.Range("I113").Value = 1679.6666
$ws.Range("J113").Value = 2238.8
$ws.Range("K113").Value = 1679.6666
$ws.Range("L113").Value = 2238.8
$ws.Range("M113").Value = 490.3334
$ws.Range("N113").Value = -6578.8

$ws.Range("H126").Value = 8971.714
$ws.Range("J126").Value = 8979.6
$ws.Range("L126").Value = 26938.8
$ws.Range("N126").Value = -31878.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

$ws.Range("H126").Value = 2694.8572
$ws.Range("J126").Value = 2998
$ws.Range("L126").Value = 8994
$ws.Range("N126").Value = -13934

$ws.Range("H132").Value = 2521.7817
$ws.Range("I132").Value = 2621.75
$ws.Range("J132").Value = 2332.3684
$ws.Range("K132").Value = 7865.25
$ws.Range("L132").Value = 6997.1052
$ws.Range("M132").Value = -5335.25
$ws.Range("N132").Value = -12057.1052
